$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 40. This shifts old rows 40..152 down to 41..153,
# matching the rest of the row's fixed fields (market/category metadata), while
# the new row 40 receives fresh observation data (date, volume, prices).
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 using the template of the row now sitting
# at 41 (the old row 40, shifted down) for the fields that stay constant, and
# the new values for the fields that actually changed.
$ws.Range("A40").Value = $ws.Range("A41").Value2
$ws.Range("B40").Value = $ws.Range("B41").Value2
$ws.Range("C40").Value = $ws.Range("C41").Value2
$ws.Range("D40").Value = 44622
$ws.Range("E40").Value = $ws.Range("E41").Value2
$ws.Range("F40").Value = $ws.Range("F41").Value2
$ws.Range("G40").Value = $ws.Range("G41").Value2
$ws.Range("H40").Value = $ws.Range("H41").Value2
$ws.Range("I40").Value = $ws.Range("I41").Value2
$ws.Range("J40").Value = 35
$ws.Range("K40").Value = 22000
$ws.Range("L40").Value = 22000
$ws.Range("M40").Value = 22000
$ws.Range("N40").Value = $ws.Range("N41").Value2
$ws.Range("O40").Value = $ws.Range("O41").Value2
$ws.Range("P40").Value = 880
$ws.Range("Q40").Value = $ws.Range("Q41").Value2
$ws.Range("R40").Value = $ws.Range("R41").Value2

# Match the date cell's number format/style (column D uses style index 2 for dates).
$ws.Range("D40").NumberFormat = $ws.Range("D41").NumberFormat
